$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.083.13"
$ws.Range("E2").Value = "  +4.50%  "
$ws.Range("D3").Value = "2.361.37"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.88%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "2.358.66"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("E10").Value = "  +8.05%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.72%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.783.58"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "57.021.60"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "2.369.24"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.11%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  +6.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "0.0₃0739"
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("E30").Value = "  +9.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.22%  "
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("E37").Value = "  +4.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.03%  "
$ws.Range("E39").Value = "  +7.20%  "
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.51%  "
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "277.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.65%  "
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.382"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.62%  "
